# Applies the cryptos-list price/volume/coin-ranking refresh described by the commit
# "Updated cryptos list on Tue Jan 30 04:26:09 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.489.59'
$ws.Range('E2').Value = '  +2.89%  '

# Row 3
$ws.Range('D3').Value = '2.310.44'
$ws.Range('E3').Value = '  +1.82%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').Value = '''310.89'
$ws.Range('E5').Value = '  +1.59%  '

# Row 6
$ws.Range('D6').Value = '''101.88'

# Row 7
$ws.Range('D7').Value = '''0.537'
$ws.Range('E7').Value = '  +1.59%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').Value = '''0.530'
$ws.Range('E9').Value = '  +7.48%  '

# Row 10
$ws.Range('D10').Value = '''35.79'
$ws.Range('E10').Value = '  +1.65%  '

# Row 11
$ws.Range('D11').Value = '''0.0815'
$ws.Range('E11').Value = '  +3.15%  '

# Row 12
$ws.Range('D12').Value = '''0.112'
$ws.Range('E12').Value = '  -0.44%  '

# Row 13
$ws.Range('D13').Value = '''7.00'
$ws.Range('E13').Value = '  +0.63%  '

# Row 14
$ws.Range('D14').Value = '2.667.91'
$ws.Range('E14').Value = '  +1.77%  '

# Row 15
$ws.Range('D15').Value = '''14.98'
$ws.Range('E15').Value = '  +1.75%  '

# Row 16
$ws.Range('D16').Value = '2.309.65'
$ws.Range('E16').Value = '  +0.64%  '

# Row 17
$ws.Range('D17').Value = '''0.808'
$ws.Range('E17').Value = '  +2.02%  '

# Row 18
$ws.Range('D18').Value = '43.396.99'
$ws.Range('E18').Value = '  +2.94%  '

# Row 19
$ws.Range('D19').Value = '''12.36'
$ws.Range('E19').Value = '  -0.09%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0928'
$ws.Range('E20').Value = '  +2.44%  '

# Row 21
$ws.Range('D21').Value = '''6.17'
$ws.Range('E21').Value = '  +2.49%  '

# Row 22
$ws.Range('D22').Value = '''68.14'
$ws.Range('E22').Value = '  +0.20%  '

# Row 23
$ws.Range('D23').Value = '''241.66'

# Row 24
$ws.Range('E24').Value = '  +3.23%  '

# Row 25
$ws.Range('D25').Value = '''2.62'
$ws.Range('E25').Value = '  +1.62%  '

# Row 27
$ws.Range('E27').Value = '  -1.59%  '

# Row 28
$ws.Range('D28').Value = '''24.66'
$ws.Range('E28').Value = '  +4.51%  '

# Row 29
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '''36.74'
$ws.Range('E29').Value = '  -2.85%  '

# Row 30
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '''9.66'
$ws.Range('E30').Value = '  +1.35%  '

# Row 31
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = '''2.12'
$ws.Range('E31').Value = '  -0.04%  '

# Row 32
$ws.Range('D32').Value = '''167.58'
$ws.Range('E32').Value = '  +3.98%  '

# Row 33
$ws.Range('D33').Value = '''5.28'
$ws.Range('E33').Value = '  +0.85%  '

# Row 34
$ws.Range('E34').Value = '  +0.06%  '

# Row 35
$ws.Range('D35').Value = '''0.0745'
$ws.Range('E35').Value = '  +0.88%  '

# Row 36
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''2.50'
$ws.Range('E36').Value = '  +5.39%  '

# Row 37
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '''3.08'
$ws.Range('E37').Value = '  -2.45%  '

# Row 38
$ws.Range('D38').Value = '''17.57'
$ws.Range('E38').Value = '  -0.25%  '

# Row 39
$ws.Range('D39').Value = '''0.106'
$ws.Range('E39').Value = '  +1.30%  '

# Row 40
$ws.Range('D40').Value = '''1.87'
$ws.Range('E40').Value = '  +2.77%  '

# Row 41
$ws.Range('E41').Value = '  +1.60%  '

# Row 42
$ws.Range('D42').Value = '''4.34'
$ws.Range('E42').Value = '  +6.54%  '

# Row 43
$ws.Range('E43').Value = '  -1.02%  '

# Row 44
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0289'
$ws.Range('E44').Value = '  +2.78%  '

# Row 45
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.967.38'
$ws.Range('E45').Value = '  +1.26%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''19.22'
$ws.Range('E46').Value = '  -0.56%  '

# Row 47
$ws.Range('D47').Value = '''2.98'
$ws.Range('E47').Value = '  +2.54%  '

# Row 48
$ws.Range('D48').Value = '''10.01'
$ws.Range('E48').Value = '  +1.17%  '

# Row 49
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''55.69'
$ws.Range('E49').Value = '  +4.01%  '

# Row 50
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.92'
$ws.Range('E50').Value = '  +5.30%  '

# Row 51
$ws.Range('D51').Value = '''1.57'
$ws.Range('E51').Value = '  +7.02%  '
